$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data source corrected: columns J and K (rows 1-51) are now uniformly 0.3.
# (Column J previously held a mix of a text header "r" in J1 and 1s in J2:J51;
# column K previously held a text header "s" in K1 and 0.3 in K2:K51.)
$ws.Range("J1:K51").Value = 0.3

# Update the view/selection to match the saved state: scrolled so row 38 is
# at the top, with K1:K51 selected (active cell K1).
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K1:K51").Select()
